$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 18, shifting the existing rows 18-42 down to 19-43.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Cells.Item(18, 1).Value = 6
$ws.Cells.Item(18, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 44799
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = 100112035
$ws.Cells.Item(18, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 220
$ws.Cells.Item(18, 11).Value = 18000
$ws.Cells.Item(18, 12).Value = 20000
$ws.Cells.Item(18, 13).Value = 18909
$ws.Cells.Item(18, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(18, 16).Value = 1261
$ws.Cells.Item(18, 17).Value = 15
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Ensure date-style formatting (style used by the other "Fecha" cells) is applied.
$ws.Range("D18").NumberFormat = $ws.Range("D19").NumberFormat
